$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.616.02"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "1.600.59"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").Value = "212.21"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").Value = "26.86"
$ws.Range("E8").Value = "  +3.58%  "
$ws.Range("E9").Value = "  +1.31%  "
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").Value = "1.828.66"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "1.604.90"
$ws.Range("E13").Value = "  +2.06%  "
$ws.Range("D14").Value = "29.607.93"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "0.538"
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").Value = "63.89"
$ws.Range("E17").Value = "  +2.50%  "
$ws.Range("D18").Value = "241.53"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").Value = "7.63"
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").Value = "9.22"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").Value = "154.95"
$ws.Range("D26").Value = "15.34"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("D30").Value = "0.0478"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("D34").Value = "1.424.78"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("E36").Value = "  +4.88%  "
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("E40").Value = "  +2.92%  "
$ws.Range("D41").Value = "1.97"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  +5.20%  "
$ws.Range("D43").Value = "53.98"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("E44").Value = "  +2.36%  "
$ws.Range("D46").Value = "0.982"
$ws.Range("E46").Value = "  +15.85%  "
$ws.Range("D47").Value = "66.26"
$ws.Range("E47").Value = "  +2.81%  "
$ws.Range("D48").Value = "5.31"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("D49").Value = "1.739.43"
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("D50").Value = "86.01"
$ws.Range("E51").Value = "  +6.16%  "
